$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the punctuation-laden text values in column B (rows 2-5) with
# plain numeric values (the punctuation characters are stripped out).
$ws.Range("B2").Value = 506542
$ws.Range("B3").Value = 146271
$ws.Range("B4").Value = 919883
$ws.Range("B5").Value = 710350

# Update the active cell selection to B3.
$ws.Range("B3").Select()
